# Add explicit root to TemplateExcelFileGenerator
#
# 1) Clarify the "root folder" explanation text (row 2).
# 2) Insert a new merged row (row 3) right below it that states the
#    explicit root path, pushing the Pol0/Pol45/Pol90/Pol135 header row
#    down from row 4 to row 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the explanatory text in row 2 to point at an explicit root.
$ws.Range("A2").Value = "The path must be the remaining path after the root folder, which is:"

# Make room for the new row: shift the old row 4 (Pol0/Pol45/Pol90/Pol135)
# down to row 5, leaving row 4 blank and opening up row 3.
$ws.Rows.Item(4).Insert()

# Populate the new row 3 with the explicit root folder path and merge it
# across A:F just like the two rows above it.
$ws.Range("A3").Value = "/home/masoud/Documents/four-polar/fourPolar-io/target/test-classes/fr/fresnel/fourPolar/io/imageSet/acquisition/sample/finders/excel"
$ws.Range("B3:F3").Value = ""
$ws.Range("A3:F3").Merge()
